$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "9695-TERGH"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "0"
$ws.Range("B6").Style = "Normal"

$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = "No"
$ws.Range("E6").Value = "Yes"
$ws.Range("F6").Value = "No"
$ws.Range("G6").Value = "Yes"
$ws.Range("H6").Value = "No"
$ws.Range("I6").Value = "No"
$ws.Range("J6").Value = "No"
$ws.Range("K6").Value = "Month-to-month"
$ws.Range("L6").Value = "Yes"
$ws.Range("M6").Value = "Electronic check"
$ws.Range("N6").Value = 96.05
$ws.Range("O6").Value = 431.98
$ws.Range("P6").Value = "Betha"
$ws.Range("Q6").Value = "2024-07"
